# Add a new survey response as row 17, mirroring the formatting of the
# existing data rows (row 15, which has every column populated so its
# number formats / styles can be copied in full).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (styles + number formats) from a fully populated row
# so the new row matches the sheet's existing look (date style on A,
# general/text style on the rest, time style on L).
$ws.Range("A15:M15").Copy()
$ws.Range("A17:M17").PasteSpecial(-4122)

# Fill in the new response's values.
$ws.Range("A17").Value = 44840.866328449076
$ws.Range("B17").Value = "valc941226@gmail.com"
$ws.Range("C17").Value = "José Carlos"
$ws.Range("D17").Value = "Ocho"
$ws.Range("E17").Value = "CDMX"
$ws.Range("F17").Value = "?"
$ws.Range("G17").Value = "Hombre"
$ws.Range("H17").Value = "Clásica, Rock, Jazz"
# "Animales que me tengo" (column I) was left blank for this response.
$ws.Range("I17").Clear()
$ws.Range("J17").Value = "BCS, One Piece, Malcom"
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 0.9791666666642413
$ws.Range("M17").Value = "Si"
